# Apply updated "Return_with_prediction" (G), "return_pct_change" (H) and
# "mean_return_pct_change" (I) values to Sheet1 of the workbook, per the
# refreshed auto-recurrence computation for DAX30 annual compared returns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new G value, new H value, new I value (blank = leave as-is)
$updates = @(
    @(2,  '0.1381283694495126',  '5.103933820157454',   '7.699803064729886'),
    @(3,  '0.1481841094432054',  '66.42944253423097',   ''),
    @(4,  '-0.5920754159267158', '3.7813062809212',     ''),
    @(5,  '-0.6406194812127495', '-4.91755136554174',   ''),
    @(6,  '0.1794184707456276',  '-27.09760407487379',  ''),
    @(7,  '0.3173572464737374',  '93.72084175596636',   ''),
    @(8,  '0.1566954572801168',  '-5.210529878098106',  ''),
    @(9,  '0.2129221450538215',  '9.139050758779593',   ''),
    @(10, '-0.08114545527317552','-41.9958731130052',   ''),
    @(11, '-0.09564825617167652','19.46652630900338',   ''),
    @(12, '0.2077853665129633',  '30.65855479104192',   ''),
    @(13, '0.2455796267029775',  '19.40948179570502',   ''),
    @(14, '0.1678040005007476',  '-11.38972668673843',  ''),
    @(15, '0.2428813392719517',  '-2.810694707400112',  ''),
    @(16, '0.05413024478009123', '48.38234219339835',   ''),
    @(17, '0.02303128585570636', '-35.06934667665193',  ''),
    @(18, '0.05181027417346977', '-70.10616141954455',  ''),
    @(19, '0.1426487826520407',  '13.4357401743283',    ''),
    @(20, '0.1309168858009779',  '14.18630715587288',   ''),
    @(21, '0.1229992669726499',  '22.51066266908664',   ''),
    @(22, '0.05975216194118629', '-36.56689701654795',  ''),
    @(23, '0.08451691517183725', '-22.09703694020399',  ''),
    @(24, '-0.1766180650056217', '-41.72035297123431',  ''),
    @(25, '-0.1523574155980028', '31.51062412824701',   ''),
    @(26, '0.1897404442158593',  '19.3472605535381',    ''),
    @(27, '0.2098866037425231',  '4.699176678341075',   ''),
    @(28, '-0.01039839157509783','-29.48703126678351',  ''),
    @(29, '0.03725690099284405', '142.2820403304422',   '')
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 7).Value = [double]$u[1]
    $ws.Cells.Item($row, 8).Value = [double]$u[2]
    if ($u[3] -ne '') {
        $ws.Cells.Item($row, 9).Value = [double]$u[3]
    }
}
